$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels in row 2 ---
$ws.Range("L2").Value = "MG"
$ws.Range("M2").Value = "MIS"

# --- Column widths for J (now wider), plus new K/L columns (widths for L/M in sheet = cols 11/12) ---
$ws.Columns("J").ColumnWidth = 24.6640625
$ws.Columns("K").ColumnWidth = 25.5
$ws.Columns("L").ColumnWidth = 24.83203125

# --- Fill formulas for K3:M14 using the existing rotation pattern, creating a shared formula group ---
$ws.Range("K3:M14").Formula = "=J4"

# Row 14 wraps around to the top of each column's own data (K14 unchanged from before; L14/M14 are new)
$ws.Range("K14").Formula = "=J3"
$ws.Range("L14").Formula = "=K3"
$ws.Range("M14").Formula = "=L3"

# --- Window position tweak recorded in workbook.xml bookViews ---
$excel.Windows.Item(1).Left = 1580
$excel.Windows.Item(1).Top = 1720

# --- Selection moved past the new data (as in the diff) ---
$ws.Range("M18").Select()

Write-Output "done"
